# Presentation Cab Company Investment
# Add presenter name "Erik Perez" as a new paragraph after the date
# ("15/3/2021") in the title slide's text box, matching the existing
# run formatting (sz=2800, bold). The text box auto-fits its height, so
# the <a:ext cy="..."/> growth to accommodate the extra line happens
# automatically.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shp = $s.Shapes.Item(2)   # "TextBox 10" containing title/subtitle/date
$tr = $shp.TextFrame.TextRange

[void]$tr.InsertAfter([char]13 + "Erik Perez")
